$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "FilesTab" query (row 4, column B = "query") is updated: the `File Type`
# and `Breed` columns are removed from the Cypher RETURN clause.
$newFilesQuery = @'

MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE diag.disease_term IN ['B Cell Lymphoma']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value = $newFilesQuery

# Reflect that the edited cell is now the active selection (the workbook was
# last saved scrolled down to, and with the cursor on, the updated cell).
$ws.Range("B4").Select() | Out-Null
